$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected range to Text format first so numeric-looking strings
# (e.g. '274.85', '0.9207') are NOT auto-converted to numbers by Excel,
# matching the original inlineStr/text cell type.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = '20.248.52'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").Value = '1.442.46'
$ws.Range("E3").Value = '  +2.49%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").Value = '0.9207'
$ws.Range("E5").Value = '  -8.03%  '
$ws.Range("D6").Value = '274.85'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").Value = '0.3642'
$ws.Range("E7").Value = '  -0.97%  '
$ws.Range("D8").Value = '0.3078'
$ws.Range("E8").Value = '  -1.34%  '
$ws.Range("D9").Value = '38.69'
$ws.Range("E9").Value = '  -1.68%  '
$ws.Range("D10").Value = '1.018'
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("D11").Value = '0.06483'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").Value = '0.9985'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").Value = '5.303'
$ws.Range("E13").Value = '  -2.60%  '
$ws.Range("E14").Value = '  -1.11%  '
$ws.Range("D15").Value = '6.021'
$ws.Range("E15").Value = '  -2.33%  '
$ws.Range("D16").Value = '0.00001006'
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").Value = '1.442.30'
$ws.Range("E17").Value = '  +2.46%  '
$ws.Range("D18").Value = '0.9369'
$ws.Range("E18").Value = '  -6.37%  '
$ws.Range("D19").Value = '0.05620'
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("D20").Value = '67.47'
$ws.Range("E20").Value = '  -4.20%  '
$ws.Range("D21").Value = '5.320'
$ws.Range("E21").Value = '  -4.80%  '
$ws.Range("E22").Value = '  -3.80%  '
$ws.Range("D23").Value = '10.72'
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("E24").Value = '  -1.70%  '
$ws.Range("D25").Value = '20.273.53'
$ws.Range("E25").Value = '  +1.74%  '
$ws.Range("D26").Value = '139.17'
$ws.Range("E26").Value = '  +3.18%  '
$ws.Range("D27").Value = '2.040'
$ws.Range("E27").Value = '  -9.18%  '
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").Value = '1.594.61'
$ws.Range("E29").Value = '  +1.90%  '
$ws.Range("D30").Value = '110.11'
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("D31").Value = '4.016'
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").Value = '4.809'
$ws.Range("E32").Value = '  -9.63%  '
$ws.Range("D33").Value = '0.7792'
$ws.Range("E33").Value = '  -4.67%  '
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").Value = '1.451'
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("D36").Value = '0.05726'
$ws.Range("E36").Value = '  -1.55%  '
$ws.Range("E37").Value = '  +4.43%  '
$ws.Range("D38").Value = '4.626'
$ws.Range("E38").Value = '  -4.65%  '
$ws.Range("E39").Value = '  -4.03%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '10.13'
$ws.Range("E40").Value = '  -2.47%  '
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").Value = '0.9331'
$ws.Range("E41").Value = '  -6.74%  '
$ws.Range("D42").Value = '0.1832'
$ws.Range("E42").Value = '  -3.66%  '
$ws.Range("D43").Value = '6.932'
$ws.Range("E43").Value = '  -17.69%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.5181'
$ws.Range("E44").Value = '  -1.86%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = '3.478'
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("D46").Value = '11.75'
$ws.Range("E46").Value = '  -4.61%  '
$ws.Range("D47").Value = '115.13'
$ws.Range("E47").Value = '  +2.47%  '
$ws.Range("D48").Value = '0.5075'
$ws.Range("E48").Value = '  -1.23%  '
$ws.Range("D49").Value = '1.726'
$ws.Range("E49").Value = '  -2.31%  '
$ws.Range("E50").Value = '  +3.11%  '
$ws.Range("D51").Value = '0.9899'
$ws.Range("E51").Value = '  -1.16%  '

# Restore default cell style (the text-format trick above would otherwise
# leave a custom number format applied to these cells).
$editRange.Style = "Normal"
